$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right 6 -> 9, Wrong 3 -> 2
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12 (Total): Right 84 -> 126, Wrong -12 -> -8, Max "72/168" -> "118/252"
$ws.Range("B12").Value = 126
$ws.Range("C12").Value = -8
$ws.Range("E12").Value = "118/252"
